$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "PostFaultTime"
$ws.Range("H2").Value = "530"
$ws.Range("H2").NumberFormat = "@"

$ws.Range("H2").Select()
